$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Listado de Usuarios")
$ws2 = $wb.Worksheets.Item("Roles")

# --- 1. Populate the "Roles" sheet with the two new role rows ---
# Row 4 first (Docente, idRol 14), then row 3 (Estudiante, idRol 3) -- this
# ordering reproduces the shared-string insertion order seen in the target
# workbook (Docente ends up before Estudiante in the shared string table).
$ws2.Range("A4").Value = 14
$ws2.Range("B4").Value = "Docente"
$ws2.Range("C4").Value = "Docente"

$ws2.Range("A3").Value = 3
$ws2.Range("B3").Value = "Estudiante"
$ws2.Range("C3").Value = "Estudiante"

# Helper column D: id lookup formulas
$ws2.Range("D2").Formula = "=A2"
$ws2.Range("D3").Formula = "=A3"
$ws2.Range("D4").Formula = "=A4"

# --- 2. Defined name "roles" used by the VLOOKUP on sheet1 ---
$wb.Names.Add("roles", "=Roles!`$B`$2:`$D`$4")

# --- 3. New header + data on "Listado de Usuarios" ---
$ws1.Range("E1").Value = "Role"
$ws1.Range("E1").Font.Bold = $true
$ws1.Range("E1").Font.Color = 192
$ws1.Range("E1").Font.Size = 12
$ws1.Range("E1").Font.Name = "Calibri"

$ws1.Range("A2").Value = "Juancito"
$ws1.Range("B2").Value = "Palaviccini"
$ws1.Range("C2").Value = "pala1590"
$ws1.Range("E2").Value = "Docente"
$ws1.Range("D2").Formula = "=IFERROR(VLOOKUP(`$E2,roles,3,FALSE),"""")"

# --- 4. Column sizing for the new D/E columns ---
$ws1.Columns.Item(4).ColumnWidth = 10.1
$ws1.Columns.Item(4).Hidden = $true
$ws1.Columns.Item(5).ColumnWidth = 16.5

# --- 5. Data validation drop-down on E2:E50 referencing the Roles sheet ---
$validationRange = $ws1.Range("E2:E50")
$validationRange.Validation.Add(3, 1, 1, "=Roles!`$B`$2:`$B`$4")

# --- 6. Selection state on the Roles sheet (set while still visible) ---
$ws2.Activate()
$ws2.Range("B19").Select()

# --- 7. Hide the Roles sheet ---
$ws2.Visible = $false

# --- 8. Re-activate the main sheet and set its selection ---
$ws1.Activate()
$ws1.Range("E2").Select()
